$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Update Version and Date values on the Metadata sheet ---
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# --- Insert a new "Jurisdiction" property row (empty value) right after
#     "Contact" (row 10) and before "Description" (row 11), pushing the
#     existing Description/Purpose/Copyright/Immutable rows down by one. ---

# Capture the current rows 11-14 content before it gets overwritten.
$savedRows = @()
for ($r = 11; $r -le 14; $r++) {
    $a = $ws.Cells.Item($r, 1).Value()
    $b = $ws.Cells.Item($r, 2).Value()
    $savedRows += , @($a, $b)
}

# Write the new Jurisdiction row in the now-freed row 11 slot.
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""

# Re-write the saved rows shifted down into rows 12-15.
for ($i = 0; $i -lt $savedRows.Count; $i++) {
    $r = 12 + $i
    $ws.Cells.Item($r, 1).Value = $savedRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $savedRows[$i][1]
}

# The newly-extended row 15 needs its formatting (borders / wrap / style)
# copied over from an existing data row, since a brand new row has no
# style applied yet.
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Rename the "Include from SetOperator" sheet to "Include #0" ---
$ws2 = $wb.Worksheets.Item("Include from SetOperator")
$ws2.Name = "Include #0"

Write-Output "done"
